# Rebuild a paragraph's run structure from a set of plain-text pieces.
# Each piece becomes its own <w:r> run (no leftover <w:proofErr/> markers,
# no stray <w:rPr/> wrappers). Works by:
#   1. Inserting one brand-new, empty paragraph per piece just before the
#      existing (messy) paragraph, and giving each its text.
#   2. Merging those new paragraphs back into a single paragraph by
#      deleting the paragraph marks between them (this keeps each piece
#      as its own run instead of Word re-merging them).
#   3. Deleting the now-redundant original paragraph (and anything in it,
#      e.g. <w:proofErr/> markers) entirely.
# Returns the 1-based paragraph index of the rebuilt paragraph.
function Rebuild-Paragraph($d, $idx, $pieces) {
    $cur = $idx
    foreach ($piece in $pieces) {
        $anchor = $d.Paragraphs($cur)
        $anchor.Range.InsertParagraphBefore()
        $newp = $d.Paragraphs($cur)
        $newp.Range.Text = $piece
        $cur = $cur + 1
    }

    $firstNew = $idx
    $lastNew = $idx + $pieces.Count - 1
    for ($i = $lastNew; $i -gt $firstNew; $i--) {
        $prevPara = $d.Paragraphs($i - 1)
        $markPos = $prevPara.Range.End - 1
        $d.Range($markPos, $markPos + 1).Delete()
    }

    $mergedIdx = $firstNew
    $oldAfterMerge = $mergedIdx + 1
    $d.Paragraphs($oldAfterMerge).Range.Delete()
    return $mergedIdx
}

# Add a zero-width ("collapsed") bookmark at a character position.
# Bookmarks.Add placed directly on an already-collapsed Range lands at the
# wrong spot in this engine, so instead: insert a throwaway character,
# bookmark the (now non-empty) range around it, then delete that
# character again -- the bookmark survives, correctly collapsed in place.
function Add-CollapsedBookmark($d, $pos, $name) {
    $tmp = $d.Range($pos, $pos)
    $tmp.InsertAfter("X")
    $bmRange = $d.Range($pos, $pos + 1)
    $d.Bookmarks.Add($name, $bmRange)
    $d.Range($pos, $pos + 1).Delete()
}

$d = $word.ActiveDocument

# Work from the last affected paragraph back to the first so that earlier
# paragraph indices stay valid while later ones are being rebuilt.

# Paragraph 27: "Struktura **tablica_wskaźników – zawiera wskaźniki na
# wszystkie możliwe struktury kolejnych podruchów"
#   -> 4 runs, and the _GoBack bookmark moves here (to the very end).
$idx27 = Rebuild-Paragraph $d 27 @(
    "Struktura",
    " ",
    "– zawiera wskaźniki na wszystkie możliwe struktury kolejn",
    "ego poziomu"
)
$endPos27 = $d.Paragraphs($idx27).Range.End - 1
Add-CollapsedBookmark $d $endPos27 "_GoBack"

# Paragraph 26: "Wektor możliwe_ruchy  - współrzędne możliwych ruchów do
# wykonania na następnym poziomie" -> single run, "możliwe_ruchy" removed.
Rebuild-Paragraph $d 26 @(
    "Wektor  - współrzędne możliwych ruchów do wykonania na następnym poziomie"
) | Out-Null

# Paragraph 25: "int wynik – determinanta wyniku danego ruchu." -> 2 runs
# ("D" / "eterminanta wyniku danego ruchu."); its _GoBack bookmark is
# removed from here (it was relocated to paragraph 27 above).
Rebuild-Paragraph $d 25 @("D", "eterminanta wyniku danego ruchu.") | Out-Null

# Paragraph 24: "int x, int y - Współrzędne wykonanego ruchu" -> single run.
Rebuild-Paragraph $d 24 @("Współrzędne wykonanego ruchu") | Out-Null

# Paragraph 23: "Struktura ruchu:" -> split into 2 runs ("Struktura ruchu"
# / ":"), keeping the sz=28/szCs=28 run formatting on both.
Rebuild-Paragraph $d 23 @("Struktura ruchu", ":") | Out-Null
